# Remove time-specific parts of presentation
#
# 1) Slide 12: "within the 10-week course" -> "within the timeframe of the course"
# 2) Slide 16: remove the bullet "The tutorial will start at 12:30"

function Remove-ParagraphAt($TextRange, $Index) {
    # Total paragraph count (split on the CR paragraph separator PowerPoint
    # exposes through TextRange.Text).
    $fullText = $TextRange.Text
    $total = ($fullText.Split([char]13)).Count

    $target = $TextRange.Paragraphs($Index, 1)
    $startPos = $target.Start
    if ($Index -gt 1) {
        $prev = $TextRange.Paragraphs($Index - 1, 1)
        $startPos = $prev.Start + $prev.Length
    }

    $length = ($target.Start + $target.Length) - $startPos
    if ($Index -eq $total) {
        # The last paragraph has no trailing paragraph mark exposed through
        # .Length, so extend the deleted span by one to fully collapse it.
        $length = $length + 1
    }

    if ($length -gt 0) {
        $victim = $TextRange.Characters($startPos, $length)
        $victim.Delete()
    }
}

$p = $ppt.ActivePresentation

# --- Edit 1: slide 12 body text -------------------------------------------
$slide12 = $p.Slides.Item(12)
$body12 = $slide12.Shapes.Item(2).TextFrame.TextRange

$oldText = "within the 10-week course"
$newText = "within the timeframe of the course"
$fullText12 = $body12.Text
$idx = $fullText12.IndexOf($oldText)
if ($idx -ge 0) {
    $target12 = $body12.Characters($idx + 1, $oldText.Length)
    $target12.Text = $newText
}

# --- Edit 2: slide 16 body text -- remove "The tutorial will start at 12:30"
$slide16 = $p.Slides.Item(16)
$body16 = $slide16.Shapes.Item(2).TextFrame.TextRange

$fullText16 = $body16.Text
$marker = "The tutorial will start at 12:30"
$paragraphs16 = $fullText16.Split([char]13)
$targetIndex = -1
for ($i = 0; $i -lt $paragraphs16.Count; $i++) {
    if ($paragraphs16[$i] -eq $marker) {
        $targetIndex = $i + 1
    }
}
if ($targetIndex -gt 0) {
    Remove-ParagraphAt $body16 $targetIndex
}
